# Update the "Förändrad" (changed) date column (C) for all data rows
# from serial date 45189 to 45190 (i.e. add one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # column C ("Förändrad")
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
